$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 15.343105725750167
$ws.Range("C2").Value = 10.638188782169019
$ws.Range("D2").Value = 15.383311792664818
$ws.Range("E2").Value = 7.8425795026476139

$ws.Range("B3").Value = 13.133227646912529
$ws.Range("C3").Value = 12.959001009682119
$ws.Range("D3").Value = 12.032620589411417
$ws.Range("E3").Value = 15.496726544888144

$ws.Range("B1:E3").Select()
